$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "UserOne"
$ws.Range("B3").Value = "UserTwo"
$ws.Range("C3").Value = 35454

$ws.Range("C9").Select()
